$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '91.633.57'
$ws.Range("E2").Value = '  -2.98%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.315.28'
$ws.Range("E3").Value = '  -4.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.29'
$ws.Range("E5").Value = '  -2.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '613.97'
$ws.Range("E6").Value = '  -3.68%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.40'
$ws.Range("E7").Value = '  -2.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.383'
$ws.Range("E8").Value = '  -3.64%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.948'
$ws.Range("E10").Value = '  -2.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.316.01'
$ws.Range("E11").Value = '  -3.89%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.06'
$ws.Range("E12").Value = '  -0.03%  '

$ws.Range("E13").Value = '  -1.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.98'
$ws.Range("E14").Value = '  -3.19%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.469.82'
$ws.Range("E15").Value = '  -3.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.935.01'
$ws.Range("E16").Value = '  -4.63%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000242'
$ws.Range("E17").Value = '  -4.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.04'
$ws.Range("E18").Value = '  -3.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.313.63'
$ws.Range("E19").Value = '  -4.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.19'
$ws.Range("E20").Value = '  -2.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.83'
$ws.Range("E21").Value = '  -4.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.40'
$ws.Range("E22").Value = '  +8.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '490.29'
$ws.Range("E23").Value = '  -1.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.447'
$ws.Range("E24").Value = '  -10.74%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000182'
$ws.Range("E25").Value = '  -4.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.21'
$ws.Range("E26").Value = '  -6.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.61'
$ws.Range("E27").Value = '  -1.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.79'
$ws.Range("E28").Value = '  -1.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.494.04'
$ws.Range("E29").Value = '  -4.30%  '

$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.99'
$ws.Range("E31").Value = '  -6.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.138'
$ws.Range("E32").Value = '  +2.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.59'
$ws.Range("E33").Value = '  -5.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.55%  '

$ws.Range("E35").Value = '  -6.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.06'
$ws.Range("E36").Value = '  -7.41%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.525'
$ws.Range("E37").Value = '  -7.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '553.52'
$ws.Range("E38").Value = '  +3.01%  '

$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.31'
$ws.Range("E40").Value = '  -4.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.147'
$ws.Range("E41").Value = '  -2.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.36'
$ws.Range("E42").Value = '  -5.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.863'
$ws.Range("E43").Value = '  -7.41%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.68'
$ws.Range("E44").Value = '  -1.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.61'
$ws.Range("E45").Value = '  +3.09%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.66'
$ws.Range("E46").Value = '  -1.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0408'
$ws.Range("E47").Value = '  -0.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.38'
$ws.Range("E48").Value = '  -3.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.10'
$ws.Range("E49").Value = '  -2.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.72'
$ws.Range("E50").Value = '  -3.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.92'
$ws.Range("E51").Value = '  -0.95%  '
